$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20 (Leve Item ID 1965)
$ws.Range("H20").Value = 10523.667
$ws.Range("I20").Value = 10523.667
$ws.Range("K20").Value = 10523.667
$ws.Range("M20").Value = -10293.667

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1101.625
$ws.Range("I28").Value = 1059
$ws.Range("J28").Value = 1400
$ws.Range("K28").Value = 1059
$ws.Range("L28").Value = 1400
$ws.Range("M28").Value = -574
$ws.Range("N28").Value = -2370

# Row 35 (Leve Item ID 1965)
$ws.Range("H35").Value = 10523.667
$ws.Range("I35").Value = 10523.667
$ws.Range("K35").Value = 10523.667
$ws.Range("M35").Value = -10144.667

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3354.2727
$ws.Range("I40").Value = 2982.8333
$ws.Range("K40").Value = 2982.8333
$ws.Range("M40").Value = -2807.8333

# Row 87 (Leve Item ID 10651)
$ws.Range("H87").Value = 120354
$ws.Range("J87").Value = 120354
$ws.Range("L87").Value = 120354
$ws.Range("N87").Value = -122850

# Row 90 (Leve Item ID 10651)
$ws.Range("H90").Value = 120354
$ws.Range("J90").Value = 120354
$ws.Range("L90").Value = 361062
$ws.Range("N90").Value = -373542

# Row 96 (Leve Item ID 19894)
$ws.Range("H96").Value = 1120
$ws.Range("I96").Value = 1477
$ws.Range("J96").Value = 548.8
$ws.Range("K96").Value = 4431
$ws.Range("L96").Value = 1646.4
$ws.Range("M96").Value = -3058
$ws.Range("N96").Value = -4392.4

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1994.6
$ws.Range("I137").Value = 1951.7142
$ws.Range("J137").Value = 2094.6667
$ws.Range("K137").Value = 5855.142599999999
$ws.Range("L137").Value = 6284.000100000001
$ws.Range("M137").Value = -3305.142599999999
$ws.Range("N137").Value = -11384.0001


$ws = $wb.Worksheets.Item("ARM")
# Row 4 (Leve Item ID 5071)
$ws.Range("H4").Value = 800.5
$ws.Range("I4").Value = 1066.6666
$ws.Range("K4").Value = 1066.6666
$ws.Range("M4").Value = -950.6666

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 6285.385
$ws.Range("I61").Value = 6309.1665
$ws.Range("K61").Value = 6309.1665
$ws.Range("M61").Value = -6097.1665

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1740.6875
$ws.Range("I74").Value = 1238.1666
$ws.Range("K74").Value = 1238.1666
$ws.Range("M74").Value = -364.1666

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1740.6875
$ws.Range("I77").Value = 1238.1666
$ws.Range("K77").Value = 6190.833000000001
$ws.Range("M77").Value = -1822.833000000001

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 6285.385
$ws.Range("I136").Value = 6309.1665
$ws.Range("K136").Value = 18927.4995
$ws.Range("M136").Value = -16377.4995

# Row 139 (Leve Item ID 42321)
$ws.Range("H139").Value = 94999
$ws.Range("J139").Value = 94999
$ws.Range("L139").Value = 94999
$ws.Range("N139").Value = -105279


$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1433.5
$ws.Range("I134").Value = 1433.5
$ws.Range("K134").Value = 4300.5
$ws.Range("M134").Value = -1765.5


$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 475.44446
$ws.Range("I7").Value = 609.8333
$ws.Range("J7").Value = 206.66667
$ws.Range("K7").Value = 609.8333
$ws.Range("L7").Value = 206.66667
$ws.Range("M7").Value = -496.8333
$ws.Range("N7").Value = -432.66667

# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 62501548
$ws.Range("I16").Value = 62501548
$ws.Range("K16").Value = 62501548
$ws.Range("M16").Value = -62501261

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3127.0908
$ws.Range("I31").Value = 2362.25
$ws.Range("K31").Value = 2362.25
$ws.Range("M31").Value = -2067.25

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3127.0908
$ws.Range("I34").Value = 2362.25
$ws.Range("K34").Value = 2362.25
$ws.Range("M34").Value = -2160.25

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 815.4286
$ws.Range("I107").Value = 851.5
$ws.Range("K107").Value = 851.5
$ws.Range("M107").Value = 1068.5

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 62501548
$ws.Range("I113").Value = 62501548
$ws.Range("K113").Value = 62501548
$ws.Range("M113").Value = -62499378


$ws = $wb.Worksheets.Item("CUL")
# Row 62 (Leve Item ID 12845)
$ws.Range("H62").Value = 5652.25
$ws.Range("I62").Value = 5665
$ws.Range("J62").Value = 5614
$ws.Range("K62").Value = 16995
$ws.Range("L62").Value = 16842
$ws.Range("M62").Value = -16309
$ws.Range("N62").Value = -18214

# Row 65 (Leve Item ID 12845)
$ws.Range("H65").Value = 5652.25
$ws.Range("I65").Value = 5665
$ws.Range("J65").Value = 5614
$ws.Range("K65").Value = 50985
$ws.Range("L65").Value = 50526
$ws.Range("M65").Value = -47553
$ws.Range("N65").Value = -57390

# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 973.5
$ws.Range("I86").Value = 817.3333
$ws.Range("K86").Value = 2451.9999
$ws.Range("M86").Value = -1265.9999

# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 973.5
$ws.Range("I89").Value = 817.3333
$ws.Range("K89").Value = 7355.9997
$ws.Range("M89").Value = -1427.9997

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 2099.6667
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2099.6667
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 18897.0003
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -23957.0003

# Row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 9693.583000000001
$ws.Range("I134").Value = 1480.909
$ws.Range("K134").Value = 4442.727000000001
$ws.Range("M134").Value = 627.2729999999992


$ws = $wb.Worksheets.Item("GSM")
# Row 69 (Leve Item ID 11891)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""

# Row 72 (Leve Item ID 11891)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""

# Row 138 (Leve Item ID 42325)
$ws.Range("H138").Value = 80778
$ws.Range("J138").Value = 80778
$ws.Range("L138").Value = 80778
$ws.Range("N138").Value = -91058


$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 395.5
$ws.Range("I16").Value = 395.5
$ws.Range("K16").Value = 395.5
$ws.Range("M16").Value = -225.5

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2450
$ws.Range("I22").Value = 150
$ws.Range("J22").Value = 2910
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 2910
$ws.Range("M22").Value = 145
$ws.Range("N22").Value = -3500

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2450
$ws.Range("I27").Value = 150
$ws.Range("J27").Value = 2910
$ws.Range("K27").Value = 150
$ws.Range("L27").Value = 2910
$ws.Range("M27").Value = -43
$ws.Range("N27").Value = -3124

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 3961.3333
$ws.Range("I61").Value = 4000
$ws.Range("J61").Value = 3942
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 3942
$ws.Range("M61").Value = -3798
$ws.Range("N61").Value = -4346

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 3961.3333
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 3942
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 3942
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -8282

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2685.4285
$ws.Range("I132").Value = 2250.75
$ws.Range("J132").Value = 3265
$ws.Range("K132").Value = 6752.25
$ws.Range("L132").Value = 9795
$ws.Range("M132").Value = -4222.25
$ws.Range("N132").Value = -14855


$ws = $wb.Worksheets.Item("WVR")
# Row 46 (Leve Item ID 42037)
$ws.Range("H46").Value = 44214.5
$ws.Range("J46").Value = 44214.5
$ws.Range("L46").Value = 44214.5
$ws.Range("N46").Value = -44676.5

# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 3083.889
$ws.Range("I62").Value = 3193.1428
$ws.Range("J62").Value = 2701.5
$ws.Range("K62").Value = 3193.1428
$ws.Range("L62").Value = 2701.5
$ws.Range("M62").Value = -2569.1428
$ws.Range("N62").Value = -3949.5

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 3083.889
$ws.Range("I65").Value = 3193.1428
$ws.Range("J65").Value = 2701.5
$ws.Range("K65").Value = 15965.714
$ws.Range("L65").Value = 13507.5
$ws.Range("M65").Value = -12845.714
$ws.Range("N65").Value = -19747.5

# Row 125 (Leve Item ID 34276)
$ws.Range("H125").Value = 42916.582
$ws.Range("J125").Value = 42916.582
$ws.Range("L125").Value = 42916.582
$ws.Range("N125").Value = -52756.582

# Row 134 (Leve Item ID 42037)
$ws.Range("H134").Value = 44214.5
$ws.Range("J134").Value = 44214.5
$ws.Range("L134").Value = 132643.5
$ws.Range("N134").Value = -137713.5

